# Updated notebook, reran simulation
# ---------------------------------------------------------------------------
# This reproduces the target diff:
#  - two brand-new data rows ("Holden", "Rizzie Spiral") inserted right
#    after the "Spiral5" row (pushing "RotRing OmegaMax-90" .. "Michael-SNHex"
#    down by two rows)
#  - the running index in column A renumbered sequentially (0..29) to match
#  - "Thomas Hex" renamed to "Matthies Hex"
#  - dimension grows from A1:W29 to A1:W31
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert two new rows right before the old row 4 ("RotRing OmegaMax-90"),
#    i.e. right after row 3 ("Spiral5"). This shifts every row from the old
#    row 4 through the old row 29 down by two rows (new rows 6..31).
$ws.Rows.Item(4).Resize(2).Insert()

# Make the freshly inserted index cells (column A) look like the rest of the
# index column (bold, centered, bordered) by copying the format from A3.
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A4:A5").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# 2. Populate the new row 4 ("Holden") and row 5 ("Rizzie Spiral").
$ws.Cells.Item(4, 2).Value = "Holden"
$ws.Cells.Item(5, 2).Value = "Rizzie Spiral"

$row4vals = @(0.9678875765843284,1.011643529171346,0.9897372409060052,1.018576725654202,1.02949805566661,1.02949805566661,1.02949805566661,1.014151490799474,0.9495657332833427,0.9495657332833427,0.9913571584067444,1.02949805566661,1.014151490799474,0.9818586120414085,1.00194436585274,0.9977384265831425,0.9844848216629408,0.9977384265831425,0.9957381301638581,1.002490115264409,0.9965521888090068)
$row5vals = @(0.8676913449859891,1.04642295238159,0.9553210952298058,1.077459126798674,1.124503280208323,1.124503280208323,1.124503280208323,1.059706270955171,0.7958723791601473,0.7958723791601473,0.9635689932029645,1.124503280208323,1.059706270955171,0.9277893250576592,1.007513683092488,0.9933606434412137,0.9369665817817081,0.9933606434412137,0.9838507563883617,1.011981261152354,0.9863181803653331)

$col = 3
foreach ($val in $row4vals) {
    $ws.Cells.Item(4, $col).Value = $val
    $col++
}

$col = 3
foreach ($val in $row5vals) {
    $ws.Cells.Item(5, $col).Value = $val
    $col++
}

# 3. Renumber the running index in column A (0-based) for every data row so
#    it stays a clean 0..29 sequence across the now-30 data rows.
for ($r = 2; $r -le 31; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}

# 4. Rename "Thomas Hex" -> "Matthies Hex" (now sitting at row 11 after the
#    insert, but look it up by content so this is robust either way).
$found = $ws.Cells.Find("Thomas Hex")
if ($found -ne $null) {
    $found.Value = "Matthies Hex"
}
